$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.923.68"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "1.641.35"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.27%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "212.61"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  +0.88%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.44"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.79%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.266"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "1.873.33"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "1.640.00"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  +0.85%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.564"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "27.896.90"
$ws.Range("E17").Value = "  +1.31%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "231.83"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -0.19%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.79"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +8.81%  "
$ws.Range("E23").Value = "  +2.07%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.20%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "151.19"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("E27").Value = "  +0.27%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.69"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.18"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "1.456.02"
$ws.Range("E33").Value = "  +0.21%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.10"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  +2.95%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.563"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  +0.56%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.919"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.15%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "69.27"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -0.23%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("E46").Value = "  +6.50%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "5.34"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").Value = "1.782.93"
$ws.Range("E48").Value = "  +1.00%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "88.40"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0507"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.11%  "
